$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update existing row 3 (the VIN "1HGEM215&4" / SYMBOL_2000_CA_SELECT row)
#    MAKE_TEXT / MODEL_TEXT -> invalidVin
#    BI/PD/UM/MP SYMBOL -> I
#    ENTRYDATE -> 20010101
#    VALID -> N
# ---------------------------------------------------------------------------
$ws.Range("E3").Value2 = "invalidVin"
$ws.Range("F3").Value2 = "invalidVin"
$ws.Range("AE3").Value2 = "I"
$ws.Range("AF3").Value2 = "I"
$ws.Range("AG3").Value2 = "I"
$ws.Range("AH3").Value2 = "I"
$ws.Range("AI3").Value2 = 20010101
$ws.Range("AJ3").Value2 = "N"

# ---------------------------------------------------------------------------
# 2. Insert two brand-new rows right after row 3 (pushing the old, mostly
#    empty, row 4 down to row 6) so that the new rows 4 and 5 inherit the
#    column formatting (style) of row 3 for every column - including the
#    columns that did not previously have cells in row 4 (D,E,F,AE..AH).
# ---------------------------------------------------------------------------
$ws.Rows("4:5").Insert(-4121)

# The old (empty) row 4 has now been shifted down to row 6 and is no longer
# needed - the new rows 4 and 5 fully replace it with real data.
$ws.Rows("6").Delete()

# ---------------------------------------------------------------------------
# 3. Populate new row 4 - duplicate of the original row 3 data but with
#    ENTRYDATE = 20020101 (BI/PD/UM/MP SYMBOL stay "O", VALID stays "Y")
# ---------------------------------------------------------------------------
$ws.Range("A4").Value2 = "1HGEM215&4"
$ws.Range("B4").Value2 = "SYMBOL_2000_CA_SELECT"
$ws.Range("C4").Value2 = 2005
$ws.Range("D4").Value2 = "TEST"
$ws.Range("E4").Value2 = "TEST"
$ws.Range("F4").Value2 = "TEST"
$ws.Range("G4").Value2 = "MDX ADVANCE"
$ws.Range("H4").Value2 = 53080
$ws.Range("I4").Value2 = "WAG"
$ws.Range("J4").Value2 = "TEST"
$ws.Range("K4").Value2 = "TEST"
$ws.Range("L4").Value2 = "TEST"
$ws.Range("M4").Value2 = "WAG"
$ws.Range("N4").Value2 = "3.5L V6   "
$ws.Range("O4").Value2 = 6
$ws.Range("P4").Value2 = "G"
$ws.Range("Q4").Value2 = 214
$ws.Range("R4").Value2 = "2WD"
$ws.Range("S4").Value2 = 2
$ws.Range("T4").Value2 = "000R"
$ws.Range("U4").Value2 = "DUAL AIR BAGS FRONT"
$ws.Range("V4").Value2 = 2
$ws.Range("W4").Value2 = "4 WHEEL STANDARD"
$ws.Range("X4").Value2 = "STD"
$ws.Range("Y4").Value2 = "B-IMMOBILIZER/KEYLSS ENTRY/ALARM"
$ws.Range("Z4").Value2 = "I"
$ws.Range("AA4").Value2 = 39
$ws.Range("AB4").Value2 = 43
$ws.Range("AC4").Value2 = "A"
$ws.Range("AD4").Value2 = "Y"
$ws.Range("AE4").Value2 = "O"
$ws.Range("AF4").Value2 = "O"
$ws.Range("AG4").Value2 = "O"
$ws.Range("AH4").Value2 = "O"
$ws.Range("AI4").Value2 = 20020101
$ws.Range("AJ4").Value2 = "Y"
$ws.Range("AK4").Value2 = "N"
$ws.Range("AL4").Value2 = "N"

# ---------------------------------------------------------------------------
# 4. Populate new row 5 - another duplicate of row 3, "secondValid" VIN row,
#    with ENTRYDATE = 20150101, BI/PD/UM/MP SYMBOL = "S", VALID stays "Y"
# ---------------------------------------------------------------------------
$ws.Range("A5").Value2 = "1HGEM215&4"
$ws.Range("B5").Value2 = "SYMBOL_2000_CA_SELECT"
$ws.Range("C5").Value2 = 2005
$ws.Range("D5").Value2 = "TEST"
$ws.Range("E5").Value2 = "secondValid"
$ws.Range("F5").Value2 = "secondValid"
$ws.Range("G5").Value2 = "MDX ADVANCE"
$ws.Range("H5").Value2 = 53080
$ws.Range("I5").Value2 = "WAG"
$ws.Range("J5").Value2 = "TEST"
$ws.Range("K5").Value2 = "TEST"
$ws.Range("L5").Value2 = "TEST"
$ws.Range("M5").Value2 = "WAG"
$ws.Range("N5").Value2 = "3.5L V6   "
$ws.Range("O5").Value2 = 6
$ws.Range("P5").Value2 = "G"
$ws.Range("Q5").Value2 = 214
$ws.Range("R5").Value2 = "2WD"
$ws.Range("S5").Value2 = 2
$ws.Range("T5").Value2 = "000R"
$ws.Range("U5").Value2 = "DUAL AIR BAGS FRONT"
$ws.Range("V5").Value2 = 2
$ws.Range("W5").Value2 = "4 WHEEL STANDARD"
$ws.Range("X5").Value2 = "STD"
$ws.Range("Y5").Value2 = "B-IMMOBILIZER/KEYLSS ENTRY/ALARM"
$ws.Range("Z5").Value2 = "I"
$ws.Range("AA5").Value2 = 39
$ws.Range("AB5").Value2 = 43
$ws.Range("AC5").Value2 = "A"
$ws.Range("AD5").Value2 = "Y"
$ws.Range("AE5").Value2 = "S"
$ws.Range("AF5").Value2 = "S"
$ws.Range("AG5").Value2 = "S"
$ws.Range("AH5").Value2 = "S"
$ws.Range("AI5").Value2 = 20150101
$ws.Range("AJ5").Value2 = "Y"
$ws.Range("AK5").Value2 = "N"
$ws.Range("AL5").Value2 = "N"

# ---------------------------------------------------------------------------
# 5. Column F width - a new custom width is introduced for this column.
# ---------------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 9.8

# ---------------------------------------------------------------------------
# 6. Update selection to match the authored state.
# ---------------------------------------------------------------------------
$ws.Range("N16").Select()
